# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Rewrites the worker mora table (B16:G23 on "Hoja1") so that it reflects
# the updated dataset: existing workers get new/adjusted "Periodo Mora" /
# "Valor Mora" / "Salario Basico" rows, and two new workers
# (ANTONIO JOSE TOBAR FUENTES / JHON ALEXANDER GIL FRANCO) are folded
# into the visible range while the period code for KEVIN / SIXTO / CATALINA
# moves from 2507 to 2508.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Columns: B=Tipo Doc, C=N Doc Trabajador, D=Nombre Trabajador,
#          E=Periodo Mora, F=Valor Mora, G=Salario Basico
# (B stays "CC" for every row - unchanged by this edit)

$ws.Range("C16").Value = "1047424489"
$ws.Range("D16").Value = "CARLOS EDUARDO MENA VEGA"
$ws.Range("E16").Value = "1609"
$ws.Range("F16").Value = 160000
$ws.Range("G16").Value = 3000000

$ws.Range("C17").Value = "1047424489"
$ws.Range("D17").Value = "CARLOS EDUARDO MENA VEGA"
$ws.Range("E17").Value = "1609"
$ws.Range("F17").Value = 6435
$ws.Range("G17").Value = 689454

$ws.Range("C18").Value = "1047457676"
$ws.Range("D18").Value = "JHON ALEXANDER GIL FRANCO"
$ws.Range("E18").Value = "1807"
$ws.Range("F18").Value = 18638
$ws.Range("G18").Value = 1164900

$ws.Range("C19").Value = "1050964767"
$ws.Range("D19").Value = "DAVID ENRIQUE LOZANO AGUAS"
$ws.Range("E19").Value = "2212"
$ws.Range("F19").Value = 80400
$ws.Range("G19").Value = 2010000

$ws.Range("C20").Value = "1050964767"
$ws.Range("D20").Value = "DAVID ENRIQUE LOZANO AGUAS"
$ws.Range("E20").Value = "2301"
$ws.Range("F20").Value = 80400
$ws.Range("G20").Value = 2010000

$ws.Range("C21").Value = "1143354710"
$ws.Range("D21").Value = "KEVIN NARVAEZ OSPINO"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 56940
$ws.Range("G21").Value = 1423500

$ws.Range("C22").Value = "8850805"
$ws.Range("D22").Value = "SIXTO HUMBERTO PERIÑAN SILVA"
$ws.Range("E22").Value = "2508"
$ws.Range("F22").Value = 56940
$ws.Range("G22").Value = 1423500

$ws.Range("C23").Value = "23071454"
$ws.Range("D23").Value = "CATALINA DEL CARMEN SALAS ELLES"
$ws.Range("E23").Value = "2508"
$ws.Range("F23").Value = 56940
$ws.Range("G23").Value = 1423500
